$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title: "Version 2" -> "Version 3"   (paragraph 3)
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute("Version 2", $true, $false, $false, $false, $false, $true, 1, $false, "Version 3", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "...convert the logical schema..." -> "...convert the optimized
#    logical schema..."   (paragraph 5)
# ------------------------------------------------------------------
$ins = $d.Paragraphs.Item(5).Range
$ins.Find.Execute("logical schema into a physical schema") | Out-Null
$ins.Collapse(1)
$ins.InsertBefore("optimized ")

# The '_GoBack' bookmark follows the newly typed text, right before
# "logical schema..." -- adding it here relocates it from wherever it
# used to be (Word only ever keeps a single '_GoBack').
$gb = $d.Paragraphs.Item(5).Range
$gb.Find.Execute("logical schema into a physical schema") | Out-Null
$gb.Collapse(1)
$d.Bookmarks.Add("_GoBack", $gb) | Out-Null

# ------------------------------------------------------------------
# 3. "Based on: Logical Schema(Version 3)" ->
#    "Based on: Optimized Logical Schema(Version 4)"   (paragraph 7)
# ------------------------------------------------------------------
$p7ins = $d.Paragraphs.Item(7).Range
$p7ins.Find.Execute("Logical Schema(Version") | Out-Null
$p7ins.Collapse(1)
$p7ins.InsertBefore("Optimized ")

$p7ver = $d.Paragraphs.Item(7).Range
$p7ver.Find.Execute("Version 3", $true, $false, $false, $false, $false, $true, 1, $false, "Version 4", 2) | Out-Null

# ------------------------------------------------------------------
# 4. "Changes: Two derived tables, COrder and ForPart were added." ->
#    "Changes: No changes were made to the Optimized Logical Schema"
#    (paragraph 9)
# ------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9).Range
$p9.Find.Execute("Two derived tables, COrder and ForPart were added.", $true, $false, $false, $false, $false, $true, 1, $false, "No changes were made to the Optimized Logical Schema", 2) | Out-Null
